$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42647.68141203704
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 10039.33
$ws.Range("D6").Value = 10085.219999999999
$ws.Range("E6").Value = 104.06
$ws.Range("F6").Value = 103.11
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = -0.91
$ws.Range("I6").Value = $false
